{"js": "// Split the Title, Author and Abstract paragraph text into one run per\n// word, with the separating spaces as their own runs too (text content is\n// unchanged \u2014 only the run boundaries change).\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build a run-per-token OOXML fragment: words become their own <w:r>, and\n// the single space between words becomes its own <w:r> too.\nfunction buildRunsOoxml(text) {\n  const tokens = text.split(/( )/).filter((t) => t.length > 0);\n  return tokens\n    .map(\n      (t) =>\n        `<w:r><w:t xml:space=\"preserve\">${escapeXml(t)}</w:t></w:r>`\n    )\n    .join(\"\");\n}\n\nfunction wrapParagraphOoxml(pStyleVal, innerRunsXml) {\n  const pPr = pStyleVal\n    ? `<w:pPr><w:pStyle w:val=\"${escapeXml(pStyleVal)}\"/></w:pPr>`\n    : \"\";\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    `<w:p>${pPr}${innerRunsXml}</w:p>` +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Map each target paragraph (identified by its current, unsplit text) to\n// its replacement plain text (identical content, just re-run below).\nconst targets = [\n  {\n    match: \"Questions: Introduction to quadratic equations\",\n    text: \"Questions: Introduction to quadratic equations\",\n  },\n  {\n    match: \"Tom Coleman\",\n    text: \"Tom Coleman\",\n  },\n  {\n    match:\n      \"A selection of questions for the study guide on introduction to quadratic equations.\",\n    text:\n      \"A selection of questions for the study guide on introduction to quadratic equations.\",\n  },\n];\n\nfor (const para of paragraphs.items) {\n  const target = targets.find((t) => t.match === para.text);\n  if (!target) continue;\n\n  para.load(\"style\");\n  await context.sync();\n  const styleName = para.style;\n\n  const runsXml = buildRunsOoxml(target.text);\n  const ooxml = wrapParagraphOoxml(styleName, runsXml);\n\n  para.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Split the Title, Author and Abstract paragraph text into one run per\n# word, with the separating spaces as their own runs too (the text\n# content itself is unchanged -- only the run boundaries change).\n\n$d = $word.ActiveDocument\n\nfunction Escape-Xml($s) {\n  $s = $s -replace '&', '&amp;'\n  $s = $s -replace '<', '&lt;'\n  $s = $s -replace '>', '&gt;'\n  $s = $s -replace '\"', '&quot;'\n  return $s\n}\n\nfunction Build-RunsXml($text) {\n  # NB: the engine's -split '( )' does not keep the captured delimiter\n  # (unlike real PowerShell), so split on plain spaces and re-insert a\n  # dedicated space run between each pair of words ourselves.\n  # NB 2: this interpreter leaks loop-variable names across function\n  # scopes (no per-function local scope for the `for` counter), so this\n  # inner loop must NOT reuse \"$i\" -- the caller's outer loop uses that\n  # name and would otherwise get its counter clobbered.\n  $words = $text -split ' '\n  $xml = \"\"\n  for ($wordIdx = 0; $wordIdx -lt $words.Count; $wordIdx++) {\n    if ($wordIdx -gt 0) {\n      $xml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n    }\n    $xml += '<w:r><w:t xml:space=\"preserve\">' + (Escape-Xml $words[$wordIdx]) + '</w:t></w:r>'\n  }\n  return $xml\n}\n\nfunction Replace-ParagraphWithRuns($para, $newText) {\n  $styleId = $para.Range.ParagraphStyle.NameLocal\n  $pPrXml = \"\"\n  if ($styleId) {\n    $pPrXml = '<w:pPr><w:pStyle w:val=\"' + (Escape-Xml $styleId) + '\"/></w:pPr>'\n  }\n  $runsXml = Build-RunsXml $newText\n  $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + $pPrXml + $runsXml + '</w:p>'\n  $null = $para.Range.InsertXML($xml)\n}\n\n$targets = @{\n  \"Questions: Introduction to quadratic equations\" = \"Questions: Introduction to quadratic equations\"\n  \"Tom Coleman\" = \"Tom Coleman\"\n  \"A selection of questions for the study guide on introduction to quadratic equations.\" = \"A selection of questions for the study guide on introduction to quadratic equations.\"\n}\n\n# Paragraph.Range.InsertXML replaces a paragraph's content in place, so\n# walking forward by index is safe (paragraph count/order is unaffected).\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $para = $d.Paragraphs($i)\n  $text = $para.Range.Text\n  # Paragraph.Range.Text includes the trailing paragraph mark; strip it.\n  $trimmed = $text.TrimEnd(\"`r\", \"`a\")\n  if ($targets.ContainsKey($trimmed)) {\n    Replace-ParagraphWithRuns $para $targets[$trimmed]\n  }\n}\n"}
